# BASE - VENDA E REGISTRO.xlsx
# ---------------------------------------------------------------------------
# The four "CONTRATAÇÃO / PREV." dates in AF17:AF20 are replaced with a
# literal "-" placeholder (a quoted dash, stored as text). The cells
# downstream (AJ/AK/AN) are driven by shared formulas such as
#   AJ: =IFERROR(WORKDAY(AF#,1),"-")
#   AK: =IFERROR(WORKDAY(AJ#,($AM$2-1)),"-")
#   AN: =IFERROR(AM#-AK#,"-")
# so once AF# stops being a date, those formulas naturally fall through to
# their "-" fallback on recalculation - nothing else needs to be touched by
# hand.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF17").Value = '"-"'
$ws.Range("AF18").Value = '"-"'
$ws.Range("AF19").Value = '"-"'
$ws.Range("AF20").Value = '"-"'

# ---------------------------------------------------------------------------
# Move the view: the worksheet had scrolled to column S with AG27 selected;
# the saved state now shows column P with AJ20 selected.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AJ20").Select()
